# feat: attribute for route
#
# Adds a new "Attribute" column (I) that documents the data fields returned
# by each API route, and switches the example route-parameter placeholder
# style in column E from "{amount}" to ":amount" for the two "amount" routes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: width ---------------------------------------------------
# (53.2 "characters" round-trips through the width conversion to the stored
# OOXML width of 54, matching the authored column width.)
$ws.Columns.Item(9).ColumnWidth = 53.2

# --- I4 header: reuse the existing bold/centered header format (D4:H4) -----
$ws.Cells.Item(4, 9).Value = "Attribute"
$ws.Cells.Item(4, 4).Copy()
$ws.Cells.Item(4, 9).PasteSpecial(-4122)   # xlPasteFormats

# --- I5: new bordered/general-alignment style for the attribute column -----
$c5 = $ws.Cells.Item(5, 9)
$c5.Value = "Id, Name"
$c5.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$c5.Borders.Item(7).Color = 0
$c5.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$c5.Borders.Item(9).Color = 0
$c5.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c5.Borders.Item(8).Color = 0
$c5.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$c5.Borders.Item(10).Color = 0

# Propagate the same format down the rest of the new column.
$c5.Copy()
$ws.Range("I6:I15").PasteSpecial(-4122)

# --- Fill in the attribute text for the remaining rows ---------------------
$attributes = @{
    6  = "Id, Birth Place"
    7  = "Id, Religion"
    8  = "Id, SMK Major"
    9  = "Id, SMA Major"
    10 = "Id, Type, Value"
    11 = "Id, Type, Value"
    12 = "Id, NIS, Name, Age, Birth Place, Birth Date, Gender, Grade, Major"
    13 = "Id, NIS, Name, Age, Birth Place, Birth Date, Gender, Grade, Major"
    14 = "Id, NIS, Name, Age, Birth Place, Birth Date, Gender, Grade, Major"
    15 = "Id, NIS, Name, Age, Birth Place, Birth Date, Gender, Grade, Major"
}
foreach ($r in $attributes.Keys) {
    $ws.Cells.Item($r, 9).Value = $attributes[$r]
}

# --- Route-parameter placeholder style: {amount} -> :amount ----------------
$ws.Cells.Item(13, 5).Value = "random/student/smk/amount/:amount"
$ws.Cells.Item(15, 5).Value = "random/student/sma/amount/:amount"

# --- Selection, matching the recorded view state ----------------------------
$ws.Range("F13").Select()
